# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price record is inserted at the top of the data table
# (row 11, right after the header row). This pushes the existing rows
# 11-110 down by one (to rows 12-111). The new record reuses the same
# product/market attributes as the row that used to be first, but with
# an updated date (Fecha) and origin (Origen).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 11:110 down to 12:111, duplicating row 11's
# formatting (incl. the date style on column D) into the freshly
# inserted blank row.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly entry.
$ws.Range("A11").Value2 = 5
$ws.Range("B11").Value2 = "Macroferia Regional de Talca"
$ws.Range("C11").Value2 = "Maule"
$ws.Range("D11").Value2 = 44881
$ws.Range("E11").Value2 = 7
$ws.Range("F11").Value2 = 100112026
$ws.Range("G11").Value2 = "Haba"
$ws.Range("H11").Value2 = "Sin especificar"
$ws.Range("I11").Value2 = "Primera"
$ws.Range("J11").Value2 = 300
$ws.Range("K11").Value2 = 9000
$ws.Range("L11").Value2 = 9000
$ws.Range("M11").Value2 = 9000
$ws.Range("N11").Value2 = "`$/saco 25 kilos"
$ws.Range("O11").Value2 = "Región del Maule"
$ws.Range("P11").Value2 = 360
$ws.Range("Q11").Value2 = 25
$ws.Range("R11").Value2 = "Hortaliza"
